# Fix for bug in how Previously Visited positions were marked.
#
# A new "Sheet3" is added to the workbook (after the existing Sheet2) that
# contains a small worked example grid used to reason about / verify the
# "previously visited" marking logic. It becomes the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook (after the last existing sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Populate the example grid
$ws3.Range("A1").Value = 1
$ws3.Range("B1").Value = 2
$ws3.Range("C1").Value = 3
$ws3.Range("F1").Value = 5
$ws3.Range("G1").Value = 8

$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = 1
$ws3.Range("C2").Value = 9
$ws3.Range("F2").Value = 6
$ws3.Range("G2").Value = 4

$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = 6
$ws3.Range("C3").Value = 1

# Match the saved selection / make Sheet3 the active (selected) tab
$ws3.Range("B3").Select()
$ws3.Activate()
